# edit.ps1 - Apply the "That Day" revision described by the diff:
#   1. Text edit: "Nicole was cursing. That meant she was plenty angry." -> "Nicole was cursing?"
#   2. Text edit: replace "Nikki" with "her sibling" in "...in the face of Nikki's concern."
#      and delete the following sentence "Claire thought she could hear a sob in the background."
#   3. Reposition the "_GoBack" bookmark to sit right after the edited sentence (where the
#      author's cursor was left following the edit), matching real Word's behaviour of
#      tracking the most recent edit location.
#
# Only the ActiveDocument ($d) is mutated; $word / $app are pre-seeded by the harness.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Nicole was cursing. That meant she was plenty angry." -> "Nicole was cursing?"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(". That meant she was plenty angry.", $true, $false, $false, $false, $false, $true, 1, $false, "?", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "...in the face of Nikki's concern. Claire thought she could hear a
# sob in the background. Okay..." -> "...in the face of her sibling's concern. Okay..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Nikki’s concern. Claire thought she could hear a sob in the background.", $true, $false, $false, $false, $false, $true, 1, $false, "her sibling’s concern.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: move the _GoBack bookmark so it sits right before "Okay now Dad ..."
# (i.e. right after "...her sibling's concern. "), matching where the edit left off.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Okay now Dad would really disapprove.") | Out-Null
$target = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null

Write-Host "Edits applied."
